# Insert a new data row at row 39 (pushing the existing rows 39-143 down to 40-144)
# and populate it with the new weekly record described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Cells.Item(39, 1).Value2 = 10
$ws.Cells.Item(39, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value2 = "La Araucanía"
$ws.Cells.Item(39, 4).Value2 = 44414
$ws.Cells.Item(39, 5).Value2 = 9
$ws.Cells.Item(39, 6).Value2 = 100112017
$ws.Cells.Item(39, 7).Value2 = "Apio"
$ws.Cells.Item(39, 8).Value2 = "Americana (o)"
$ws.Cells.Item(39, 9).Value2 = "Primera"
$ws.Cells.Item(39, 10).Value2 = 80
$ws.Cells.Item(39, 11).Value2 = 9000
$ws.Cells.Item(39, 12).Value2 = 9000
$ws.Cells.Item(39, 13).Value2 = 9000
$ws.Cells.Item(39, 14).Value2 = "`$/docena de paquetes"
$ws.Cells.Item(39, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(39, 16).Value2 = 750
$ws.Cells.Item(39, 17).Value2 = 12
$ws.Cells.Item(39, 18).Value2 = "Hortaliza"
